$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New hour-log entries (rows 10-14) describing the April 2023 deliverables:
# update function fix, error handling, array-returning functions,
# improved/compound search, and the refreshed deliverables package.

$ws.Range("A10").Value = "Updated update function"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 45034
$ws.Range("D10").Value = "Updated the update function so it returns proper JSON data"

$ws.Range("A11").Value = "Error handling"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 45035
$ws.Range("D11").Value = "Added errors to several functions incase a facility was not found"

$ws.Range("A12").Value = "Updated functions"
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = 45037
$ws.Range("D12").Value = "Updated functions so they return tags as a array instead of string"

$ws.Range("A13").Value = "Improve search method"
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = 45037
$ws.Range("D13").Value = "Implemented compount search function properly"

$ws.Range("A14").Value = "Prepare deliverables"
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 45039
$ws.Range("D14").Value = "Included new  hour log, DB dump with sample data and postman collection"
